# daily auto push: 2026-01-23 06:52 UTC
# Insert a new data row at row 709 (shifting existing rows 709:750 down to
# 710:751), carrying the new observation for 2026/01/23 (金, hour 13, rank 114).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Push rows 709:750 down one row, opening up a blank row 709.
$ws.Rows.Item(709).Insert()

# Columns A/B hold plain text ("2026/01/23", "金") in the source data - force
# text formatting first so Excel doesn't auto-coerce the date-looking string
# into a date serial value.
$ws.Range("A709:B709").NumberFormat = "@"

$ws.Range("A709").Value = "2026/01/23"
$ws.Range("B709").Value = "金"
$ws.Range("C709").Value = 13
$ws.Range("D709").Value = 114

# Restore the default (unstyled) look so the new row matches its neighbours.
$ws.Range("A709:B709").Style = "Normal"
